$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$freq = @(562.5, 1125, 1687.5, 2718.75, 2250, 3328.125, 3000, 3562.5, 3890.625, 3187.5)
$mag  = @(1633.304321289062, 47.20366668701172, 29.50435256958008, 24.51774406433105, 17.87451934814453, 11.60765075683594, 9.267233848571777, 9.189947128295898, 9.049626350402832, 8.95585823059082)

for ($i = 0; $i -lt $freq.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $freq[$i]
    $ws.Cells.Item($row, 2).Value = $mag[$i]
}
